# Felipe Neto 1 - add "B" column (classification labels) to the Treinamento sheet
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Treinamento")

# Values for B2:B59 as introduced by the edit (integer labels for the training rows)
$values = @(3,2,4,1,2,3,4,2,4,3,3,3,4,3,3,2,1,4,4,2,2,3,2,1,4,4,4,4,4,4,2,0,4,0,0,1,4,3,1,1,2,0,3,2,4,2,3,1,3,0,4,3,1,4,2,4,4,4)

$row = 2
foreach ($v in $values) {
    $ws1.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# Update the sheet view so the scrolled position / selection matches the saved state
$ws1.Application.Goto($ws1.Range("A33"), $true)
$ws1.Range("B60").Select()
